# Remove the decorative "────...────" separator paragraphs and the
# thin empty spacer paragraphs (w:spacing w:before="40") that precede
# each "Heading3" subsection / follow each code table, as part of a
# document-wide cleanup pass.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$ranges = New-Object System.Collections.ArrayList

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    $sb = $p.SpaceBefore

    # Empty "spacer" paragraph inserted right after a code table
    # (w:pPr/w:spacing w:before="40", no text besides the paragraph mark).
    $isSpacer = ($t.Length -le 1) -and ($sb -eq 2)

    # Decorative horizontal-rule separator paragraph (long dash run).
    $isSeparator = $t -match "─────"

    if ($isSpacer -or $isSeparator) {
        [void]$ranges.Add($p.Range)
    }
}

# Delete from the end of the document backwards so earlier ranges keep
# their original Start/End offsets valid.
for ($i = $ranges.Count - 1; $i -ge 0; $i--) {
    $r = $ranges[$i]
    $r.Delete()
}
